$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A ("Match ID") - shifts all existing columns one to the right.
$ws.Columns("A").Insert()

# Header for the new column, bold (no border) like the rest of the new column's cells.
$ws.Range("A1").Value = "Match ID"
$ws.Range("A1:A19").Font.Bold = $true

# Fill the new column with the match id (29) for the visible data rows.
$ws.Range("A4:A19").Value = 29

# Rows 2, 3 and 20 are hidden; temporarily unhide them so writing a value doesn't
# stamp a custom row height, then restore their hidden state.
$ws.Rows(2).Hidden = $false
$ws.Rows(3).Hidden = $false
$ws.Rows(20).Hidden = $false

$ws.Range("A20").Value = 29

$ws.Rows(2).Hidden = $true
$ws.Rows(3).Hidden = $true
$ws.Rows(20).Hidden = $true

# Update the saved selection to match the new layout.
$ws.Range("A1:A19").Select()
